$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ronda 2 hoy 14/06: replace the Movil/Mensaje/Imagen assignment rows
# with a new batch of phone numbers, re-using the 3 existing message
# templates (Olivo / Azucena / Arboledas) in a new row order, and add
# one extra row (46) that was not present before.

$ws.Range("A2").Value = '52_7751441859'
$ws.Range("B2").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C2").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A3").Value = '52_7751456754'
$ws.Range("B3").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C3").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A4").Value = '52_7751537112'
$ws.Range("B4").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C4").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A5").Value = '52_7751551945'
$ws.Range("B5").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C5").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A6").Value = '52_7751605551'
$ws.Range("B6").Value = '📢Le presento la Privada *★Arboledas, Modelo Compostela★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯 INFORMACION GENERAL :
- 🟩 Terreno: 105 m²
- 🏡 Construcción: 114 m²
- 🛏️Habitaciones: 3 (2 en planta alta, con vestidor)
- 🚽Baños: 2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,695,000 MXN*
1. Pago Promedio mensual: $17,731.23 MXN
2. Ingreso mayor a los $47,094.01 MXN mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C6").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\arboledas.png'

$ws.Range("A7").Value = '52_7751609397'
$ws.Range("B7").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C7").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A8").Value = '52_7751618134'
$ws.Range("B8").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C8").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A9").Value = '52_7751674759'
$ws.Range("B9").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C9").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A10").Value = '52_7751820501'
$ws.Range("B10").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C10").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A11").Value = '52_7751855523'
$ws.Range("B11").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C11").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A12").Value = '52_7751916237'
$ws.Range("B12").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C12").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A13").Value = '52_7751928351'
$ws.Range("B13").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C13").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A14").Value = '52_7751937397'
$ws.Range("B14").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C14").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A15").Value = '52_7751980213'
$ws.Range("B15").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C15").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A16").Value = '52_7751995043'
$ws.Range("B16").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C16").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A17").Value = '52_7752018495'
$ws.Range("B17").Value = '📢Le presento la Privada *★Arboledas, Modelo Compostela★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯 INFORMACION GENERAL :
- 🟩 Terreno: 105 m²
- 🏡 Construcción: 114 m²
- 🛏️Habitaciones: 3 (2 en planta alta, con vestidor)
- 🚽Baños: 2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,695,000 MXN*
1. Pago Promedio mensual: $17,731.23 MXN
2. Ingreso mayor a los $47,094.01 MXN mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C17").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\arboledas.png'

$ws.Range("A18").Value = '52_7752018993'
$ws.Range("B18").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C18").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A19").Value = '52_7752026064'
$ws.Range("B19").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C19").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A20").Value = '52_7752054219'
$ws.Range("B20").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C20").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A21").Value = '52_7752055804'
$ws.Range("B21").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C21").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A22").Value = '52_7752055827'
$ws.Range("B22").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C22").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A23").Value = '52_7752061930'
$ws.Range("B23").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C23").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A24").Value = '52_7752212730'
$ws.Range("B24").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C24").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A25").Value = '52_7752282291'
$ws.Range("B25").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C25").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A26").Value = '52_7752356335'
$ws.Range("B26").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C26").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A27").Value = '52_7752359588'
$ws.Range("B27").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C27").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A28").Value = '52_7752361643'
$ws.Range("B28").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C28").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A29").Value = '52_7752509898'
$ws.Range("B29").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C29").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A30").Value = '52_7752534040'
$ws.Range("B30").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C30").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A31").Value = '52_7757067426'
$ws.Range("B31").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C31").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A32").Value = '52_7751458051'
$ws.Range("B32").Value = '📢Le presento la Privada *★Arboledas, Modelo Compostela★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯 INFORMACION GENERAL :
- 🟩 Terreno: 105 m²
- 🏡 Construcción: 114 m²
- 🛏️Habitaciones: 3 (2 en planta alta, con vestidor)
- 🚽Baños: 2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,695,000 MXN*
1. Pago Promedio mensual: $17,731.23 MXN
2. Ingreso mayor a los $47,094.01 MXN mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C32").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\arboledas.png'

$ws.Range("A33").Value = '52_7751357889'
$ws.Range("B33").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C33").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A34").Value = '52_7757548955'
$ws.Range("B34").Value = '📢Le presento la Privada *★Arboledas, Modelo Compostela★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯 INFORMACION GENERAL :
- 🟩 Terreno: 105 m²
- 🏡 Construcción: 114 m²
- 🛏️Habitaciones: 3 (2 en planta alta, con vestidor)
- 🚽Baños: 2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,695,000 MXN*
1. Pago Promedio mensual: $17,731.23 MXN
2. Ingreso mayor a los $47,094.01 MXN mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C34").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\arboledas.png'

$ws.Range("A35").Value = '52_5579934304'
$ws.Range("B35").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C35").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A36").Value = '52_7751056494'
$ws.Range("B36").Value = '📢Le presento la Privada *★Arboledas, Modelo Compostela★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯 INFORMACION GENERAL :
- 🟩 Terreno: 105 m²
- 🏡 Construcción: 114 m²
- 🛏️Habitaciones: 3 (2 en planta alta, con vestidor)
- 🚽Baños: 2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,695,000 MXN*
1. Pago Promedio mensual: $17,731.23 MXN
2. Ingreso mayor a los $47,094.01 MXN mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C36").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\arboledas.png'

$ws.Range("A37").Value = '52_7751323350'
$ws.Range("B37").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C37").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A38").Value = '52_7751373509'
$ws.Range("B38").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C38").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A39").Value = '52_7712029807'
$ws.Range("B39").Value = '📢Le presento la Privada *★Olivo, Modelo Lugo★*
 _Tulancingo, Hidalgo a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Terreno: 119m²
- 🏡 Construcción: 152 m²
- 🛏️Habitaciones: 4 (Recámara principal con walk in closet)
- 🚽Baños: 3 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
Precio *$2,400,000 MXN*
1. Pago promedio mensual $24,317.26 MXN
2. Ingreso mayor a los $59,000.00 MXN mensuales comprobables por pareja

💡¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C39").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\olivo.png'

$ws.Range("A40").Value = '52_7712059478'
$ws.Range("B40").Value = '📢Le presento la Privada *★Arboledas, Modelo Compostela★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯 INFORMACION GENERAL :
- 🟩 Terreno: 105 m²
- 🏡 Construcción: 114 m²
- 🛏️Habitaciones: 3 (2 en planta alta, con vestidor)
- 🚽Baños: 2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,695,000 MXN*
1. Pago Promedio mensual: $17,731.23 MXN
2. Ingreso mayor a los $47,094.01 MXN mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C40").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\arboledas.png'

$ws.Range("A41").Value = '52_7751036577'
$ws.Range("B41").Value = '📢Le presento la Privada *★Arboledas, Modelo Compostela★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯 INFORMACION GENERAL :
- 🟩 Terreno: 105 m²
- 🏡 Construcción: 114 m²
- 🛏️Habitaciones: 3 (2 en planta alta, con vestidor)
- 🚽Baños: 2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,695,000 MXN*
1. Pago Promedio mensual: $17,731.23 MXN
2. Ingreso mayor a los $47,094.01 MXN mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C41").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\arboledas.png'

$ws.Range("A42").Value = '52_7751244729'
$ws.Range("B42").Value = '📢Le presento la Privada *★Arboledas, Modelo Compostela★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯 INFORMACION GENERAL :
- 🟩 Terreno: 105 m²
- 🏡 Construcción: 114 m²
- 🛏️Habitaciones: 3 (2 en planta alta, con vestidor)
- 🚽Baños: 2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,695,000 MXN*
1. Pago Promedio mensual: $17,731.23 MXN
2. Ingreso mayor a los $47,094.01 MXN mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C42").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\arboledas.png'

$ws.Range("A43").Value = '52_7751267907'
$ws.Range("B43").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C43").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A44").Value = '52_7751279804'
$ws.Range("B44").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C44").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'

$ws.Range("A45").Value = '52_7751280849'
$ws.Range("B45").Value = '📢Le presento la Privada *★Arboledas, Modelo Compostela★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯 INFORMACION GENERAL :
- 🟩 Terreno: 105 m²
- 🏡 Construcción: 114 m²
- 🛏️Habitaciones: 3 (2 en planta alta, con vestidor)
- 🚽Baños: 2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,695,000 MXN*
1. Pago Promedio mensual: $17,731.23 MXN
2. Ingreso mayor a los $47,094.01 MXN mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C45").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\arboledas.png'

$ws.Range("A46").Value = '52_7751315193'
$ws.Range("B46").Value = '📢Le presento la Privada *★Azucena, Modelo Dúplex★*
 _Tulancingo, Hidalgo  a 5 Minutos de Walmart_

🎯INFORMACION GENERAL:
- 🟩Plantas: 2
- 🏡 Construcción: 95.9 m²
- 🛏️Habitaciones: 2
- 🚽Baños: 1 + 1/2
- 🚗Estacionamiento: 2
- 📺 Sala
- 🪑Comedor
- 🌿Jardín

⛳Amenidades:
- 🌲Área Verde
- 🥩Asador
- 🛝Juegos Infantiles
- 📍Céntrico: Valle Sol, 43648 Tulancingo, Hgo.

💰Información Hipotecaria & crediticia:
*$1,195,000 MXN*
1. Pago Promedio mensual $11,000.35 MXN
2. Ingreso mayor a los $25,511.47 MXN Mensuales comprobables por pareja

💡 ¿Le interesa obtener más información de esta propiedad?'
$ws.Range("C46").Value = 'C:\Users\4to Creativo\OneDrive\Desktop\WhatsAutoA\Images\azucena.png'
